$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a drug/product table in rows 4-13 (row 3 = header,
# row 14 = totals, row 15 = footer).  Two new products need to be added
# ("EMPACOZA TRIO XR 25/5/1000  30TAB" and "GLYBOFEN 5/500MG 30 F.C.TABS.")
# while keeping the existing alphabetical ordering, which pushes the
# totals/footer rows down by two rows.

# 1) Insert two blank rows right before the current totals row (row 14),
#    i.e. after the current last product row (row 13).
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# 2) Copy the formatting (styles, fonts, borders, number formats) of the
#    last existing product row into the two new rows, reusing the same
#    style definitions instead of creating new ones.
$ws.Range("A13:N13").Copy()
$ws.Range("A14:N15").PasteSpecial(-4122)

# 3) Match the row heights used by the other product rows (they alternate
#    between 24.75 and 25.5).
$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5

# 4) Re-create the merged cells for the two new rows (B:G, H:K, L:M),
#    matching the pattern used by every other product row.
$ws.Range("B14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

# 5) Re-write the product rows (9-15) so the list stays alphabetically
#    sorted with the two new products inserted in the right spots:
#      ... DOLPHIN, EMPACOZA (new), FEROGLOBIN, GLYBOFEN (new), LEZBERG,
#      RESTOHERB, VOLTAREN, سرنجات ...
$ws.Range("B9").Value = "EMPACOZA TRIO XR 25/5/1000  30TAB"
$ws.Range("H9").Value = "0:2"
$ws.Range("L9").Value = 132
$ws.Range("N9").Value = "0:0"

$ws.Range("B10").Value = "FEROGLOBIN 30 CAPS"
$ws.Range("H10").Value = "0:1"
$ws.Range("L10").Value = 90
$ws.Range("N10").Value = "0:2"

$ws.Range("B11").Value = "GLYBOFEN 5/500MG 30 F.C.TABS."
$ws.Range("H11").Value = "0:1"
$ws.Range("L11").Value = 14.25
$ws.Range("N11").Value = "0:0"

$ws.Range("B12").Value = "LEZBERG TRIO 20/5/12.5 TAB"
$ws.Range("H12").Value = "0:2"
$ws.Range("L12").Value = 38
$ws.Range("N12").Value = "0:0"

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "RESTOHERB SYRUP"
$ws.Range("H13").Value = "6:0"
$ws.Range("L13").Value = 45
$ws.Range("N13").Value = "1:0"

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H14").Value = "7:3"
$ws.Range("L14").Value = 34
$ws.Range("N14").Value = "0:3"

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "سرنجات 3 سم"
$ws.Range("H15").Value = "-2:0"
$ws.Range("L15").Value = 6
$ws.Range("N15").Value = "3:0"

# 6) Update the grand total (now on row 16) to reflect the two new rows.
$ws.Range("K16").Value = 555.08
